$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 1819.4
$ws.Range("B12").Value = 3344.91
$ws.Range("F54").Value = 26
$ws.Range("G54").Value = 580.3200000000001
$ws.Range("B61").Value = 23438.01
$ws.Range("F107").Value = 35
$ws.Range("G107").Value = 2461.2
$ws.Range("F122").Value = 70
$ws.Range("G122").Value = 7172.2
$ws.Range("B133").Value = 195876.34
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 2881.24
$ws.Range("F158").Value = 1
$ws.Range("G158").Value = 2115.59
$ws.Range("B163").Value = 37645.32
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("F174").Value = 68
$ws.Range("G174").Value = 3027.36
$ws.Range("F175").Value = 14
$ws.Range("G175").Value = 448.28
$ws.Range("B177").Value = 11123.81
$ws.Range("F195").Value = 1
$ws.Range("G195").Value = 1209
$ws.Range("B196").Value = 31149.03
$ws.Range("F221").Value = 119
$ws.Range("G221").Value = 3692.57
$ws.Range("F224").Value = 9
$ws.Range("G224").Value = 675.1799999999999
$ws.Range("B229").Value = 10783.47
$ws.Range("F247").Value = 66
$ws.Range("G247").Value = 4548.06
$ws.Range("F248").Value = 24
$ws.Range("G248").Value = 510
$ws.Range("B251").Value = 11205.46
$ws.Range("F281").Value = 147
$ws.Range("G281").Value = 4593.75
$ws.Range("F288").Value = 23
$ws.Range("G288").Value = 1886.23
$ws.Range("F311").Value = 32
$ws.Range("G311").Value = 2741.76
$ws.Range("B313").Value = 105219.84
$ws.Range("B317").Value = 61610
$ws.Range("D317").Value = 102.71
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = 62
$ws.Range("G317").Value = 6368.02
$ws.Range("B318").Value = 57077
$ws.Range("D318").Value = 93.08
$ws.Range("E318").Value = 111.2
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 93.08
$ws.Range("F327").Value = 162
$ws.Range("G327").Value = 18505.26
$ws.Range("B383").Value = 222735.07
$ws.Range("F398").Value = 277
$ws.Range("G398").Value = 8930.48
$ws.Range("F399").Value = 273
$ws.Range("G399").Value = 6120.66
$ws.Range("B404").Value = 45388.13
$ws.Range("F417").Value = 118
$ws.Range("G417").Value = 4360.1
$ws.Range("F420").Value = 595
$ws.Range("G420").Value = 83651.05
$ws.Range("F422").Value = 89
$ws.Range("G422").Value = 12625.54
$ws.Range("B423").Value = 126429.47
$ws.Range("F447").Value = 59
$ws.Range("G447").Value = 6226.86
$ws.Range("B458").Value = 35970.84
$ws.Range("F462").Value = 248
$ws.Range("G462").Value = 12851.36
$ws.Range("F470").Value = 67
$ws.Range("G470").Value = 3175.8
$ws.Range("F476").Value = 74
$ws.Range("G476").Value = 2460.5
$ws.Range("B484").Value = 132142.76
$ws.Range("F502").Value = 531
$ws.Range("G502").Value = 6802.11
$ws.Range("F503").Value = 257
$ws.Range("G503").Value = 6759.1
$ws.Range("F510").Value = 68
$ws.Range("G510").Value = 1323.28
$ws.Range("F511").Value = 896
$ws.Range("G511").Value = 5895.68
$ws.Range("F513").Value = 911
$ws.Range("G513").Value = 5912.39
$ws.Range("B518").Value = 90850.96000000001
$ws.Range("F523").Value = 30
$ws.Range("G523").Value = 1044.3
$ws.Range("F526").Value = 105
$ws.Range("G526").Value = 3655.05
$ws.Range("B528").Value = 16201.36
$ws.Range("F570").Value = 571
$ws.Range("G570").Value = 11334.35
$ws.Range("F571").Value = 195
$ws.Range("G571").Value = 1306.5
$ws.Range("B575").Value = 33970.9
$ws.Range("F631").Value = 13
$ws.Range("G631").Value = 1350.31
$ws.Range("F648").Value = 38
$ws.Range("G648").Value = 4667.54
$ws.Range("B652").Value = 203063.45
$ws.Range("F654").Value = 84
$ws.Range("G654").Value = 10966.2
$ws.Range("F655").Value = 61
$ws.Range("G655").Value = 10859.22
$ws.Range("B661").Value = 47740.92
$ws.Range("F673").Value = 2
$ws.Range("G673").Value = 4978.68
$ws.Range("B678").Value = 30435.61
$ws.Range("F698").Value = 14
$ws.Range("G698").Value = 1632.12
$ws.Range("F700").Value = 6
$ws.Range("G700").Value = 570
$ws.Range("B705").Value = 20605.06
$ws.Range("F712").Value = 4
$ws.Range("G712").Value = 247.6
$ws.Range("F724").Value = 1
$ws.Range("G724").Value = 3215.12
$ws.Range("B729").Value = 56796.68
$ws.Range("F755").Value = 28
$ws.Range("G755").Value = 2550.24
$ws.Range("B757").Value = 2550.24
$ws.Range("F768").Value = 181
$ws.Range("G768").Value = 14762.36
$ws.Range("F770").Value = 18
$ws.Range("G770").Value = 1468.08
$ws.Range("F771").Value = 210
$ws.Range("G771").Value = 27405
$ws.Range("F783").Value = 51
$ws.Range("G783").Value = 5652.33
$ws.Range("F784").Value = 378
$ws.Range("G784").Value = 51033.78
$ws.Range("F786").Value = 456
$ws.Range("G786").Value = 55043.76
$ws.Range("B788").Value = 201984.51
$ws.Range("F834").Value = 2
$ws.Range("G834").Value = 1618.64
$ws.Range("B836").Value = 2415.86
$ws.Range("F882").Value = 364
$ws.Range("G882").Value = 11003.72
$ws.Range("F883").Value = 2887
$ws.Range("G883").Value = 470898.57
$ws.Range("F885").Value = 163
$ws.Range("G885").Value = 23577.95
$ws.Range("F886").Value = 99
$ws.Range("G886").Value = 3775.86
$ws.Range("F889").Value = 171
$ws.Range("G889").Value = 11542.5
$ws.Range("B891").Value = 559974.86
$ws.Range("B897").Value = 3494275.14
$ws.Range("B898").Value = 3494275.14
